$wb = $excel.ActiveWorkbook

# --- Metadata: update "Last Updated" timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "29 Oct 2025, 03:59 PM"

# --- Top Gainers: rows 39-50 shift up by one, row 50 gets a refreshed BIL entry ---
$wsGainers = $wb.Worksheets.Item("Top Gainers")

$gainersData = @(
    @(39, "DATAMATICS", 4.9005, 7.3152, 15.7298),
    @(40, "UTKARSHBNK", 4.8768, -5.8959, -2.6215),
    @(41, "FILATEX", 4.8689, 10.274, 26.0027),
    @(42, "HITECHGEAR", 4.8651, 2.1287, 10.9905),
    @(43, "INDOTHAI", 4.8064, 4.5349, 43.748),
    @(44, "SANDUMA", 4.593, 2.1405, 30.2813),
    @(45, "LLOYDSENT", 4.5646, 1.8339, 11.234),
    @(46, "STAR", 4.5025, 4.4319, 3.662),
    @(47, "RECLTD", 4.4992, 3.4756, 3.4062),
    @(48, "NBCC", 4.4511, 3.1605, 7.6018),
    @(49, "GPPL", 4.4154, 3.4073, 5.0497),
    @(50, "BIL", 4.3654, 9.122199999999999, -0.3203)
)

foreach ($row in $gainersData) {
    $r = $row[0]
    $wsGainers.Cells.Item($r, 2).Value = $row[1]
    $wsGainers.Cells.Item($r, 3).Value = $row[2]
    $wsGainers.Cells.Item($r, 4).Value = $row[3]
    $wsGainers.Cells.Item($r, 5).Value = $row[4]
}

# --- Top Losers: rows 28-31 shift up by one, row 31 gets a refreshed SUMMITSEC entry ---
$wsLosers = $wb.Worksheets.Item("Top Losers")

$losersData = @(
    @(28, "SMSPHARMA", -3.7339, -3.0871, 17.4387),
    @(29, "BHARATWIRE", -3.5327, 22.8336, 23.8979),
    @(30, "ABSLAMC", -3.5313, -5.9355, -1.2887),
    @(31, "SUMMITSEC", -3.4113, -1.6476, 6.0097)
)

foreach ($row in $losersData) {
    $r = $row[0]
    $wsLosers.Cells.Item($r, 2).Value = $row[1]
    $wsLosers.Cells.Item($r, 3).Value = $row[2]
    $wsLosers.Cells.Item($r, 4).Value = $row[3]
    $wsLosers.Cells.Item($r, 5).Value = $row[4]
}
